$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting existing row 56 (and below) down by one.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record's data.
$ws.Cells.Item(56, 1).Value = 4
$ws.Cells.Item(56, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(56, 3).Value = "Los Lagos"
$ws.Cells.Item(56, 4).Value = 44571
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 10
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100108
$ws.Cells.Item(56, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(56, 9).Value = 100108002
$ws.Cells.Item(56, 10).Value = "Mango"
$ws.Cells.Item(56, 11).Value = "Sin especificar"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 80
$ws.Cells.Item(56, 14).Value = 8000
$ws.Cells.Item(56, 15).Value = 8500
$ws.Cells.Item(56, 16).Value = 8250
$ws.Cells.Item(56, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(56, 18).Value = "Perú"
$ws.Cells.Item(56, 19).Value = 2062
$ws.Cells.Item(56, 20).Value = 4
